# Restore C10 ("Rules" sheet, row for rule R30 / "From" value) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
